$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) stores numbers that look numeric (e.g. "607.00",
# "0.494") as literal TEXT in the source feed (note columns like "66.201.48"
# or "3.563.60" already can only be text because of the extra "."). To keep
# every updated Price cell as text (matching the original file's inlineStr
# cells) instead of Excel silently re-interpreting them as numbers, force
# each such target cell to the Text number format before writing its value.
$textFormatCells = @(
    "D5", "D6", "D9", "D12", "D14", "D15", "D19", "D20", "D21", "D22", "D24", "D27", "D29", "D30", "D32", "D33", "D35", "D38", "D39", "D40", "D41", "D43", "D44", "D46", "D47", "D48", "D49", "D50", "D51"
)
foreach ($cellRef in $textFormatCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = "66.191.72"
$ws.Range("E2").Value = "  +1.05%  "

# Row 3
$ws.Range("D3").Value = "3.562.43"
$ws.Range("E3").Value = "  +4.68%  "

# Row 4
$ws.Range("E4").Value = "  +0.06%  "

# Row 5
$ws.Range("D5").Value = "607.00"
$ws.Range("E5").Value = "  +1.73%  "

# Row 6
$ws.Range("D6").Value = "144.48"
$ws.Range("E6").Value = "  +1.88%  "

# Row 7
$ws.Range("D7").Value = "3.560.83"
$ws.Range("E7").Value = "  +4.61%  "

# Row 8
$ws.Range("E8").Value = "  +0.15%  "

# Row 9
$ws.Range("D9").Value = "0.494"
$ws.Range("E9").Value = "  +5.33%  "

# Row 10
$ws.Range("E10").Value = "  +2.47%  "

# Row 11
$ws.Range("E11").Value = "  +0.82%  "

# Row 12
$ws.Range("D12").Value = "0.413"
$ws.Range("E12").Value = "  +2.47%  "

# Row 13
$ws.Range("D13").Value = "4.167.40"
$ws.Range("E13").Value = "  +4.85%  "

# Row 14
$ws.Range("D14").Value = "0.0000207"
$ws.Range("E14").Value = "  +4.02%  "

# Row 15
$ws.Range("D15").Value = "30.08"
$ws.Range("E15").Value = "  +2.49%  "

# Row 16
$ws.Range("D16").Value = "3.564.12"
$ws.Range("E16").Value = "  +4.93%  "

# Row 17
$ws.Range("D17").Value = "66.307.45"
$ws.Range("E17").Value = "  +1.16%  "

# Row 18
$ws.Range("E18").Value = "  -0.49%  "

# Row 19
$ws.Range("D19").Value = "11.34"
$ws.Range("E19").Value = "  +9.90%  "

# Row 20
$ws.Range("D20").Value = "6.19"
$ws.Range("E20").Value = "  +1.72%  "

# Row 21
$ws.Range("D21").Value = "14.89"
$ws.Range("E21").Value = "  +2.63%  "

# Row 22
$ws.Range("D22").Value = "431.26"
$ws.Range("E22").Value = "  +4.26%  "

# Row 23
$ws.Range("E23").Value = "  +6.27%  "

# Row 24
$ws.Range("D24").Value = "78.91"
$ws.Range("E24").Value = "  +2.61%  "

# Row 25
$ws.Range("D25").Value = "3.704.80"
$ws.Range("E25").Value = "  +4.75%  "

# Row 26
$ws.Range("E26").Value = "  -0.02%  "

# Row 27
$ws.Range("D27").Value = "0.0000118"
$ws.Range("E27").Value = "  +9.40%  "

# Row 28
$ws.Range("E28").Value = "  +4.90%  "

# Row 29
$ws.Range("D29").Value = "7.99"
$ws.Range("E29").Value = "  +3.13%  "

# Row 30
$ws.Range("D30").Value = "9.10"
$ws.Range("E30").Value = "  -0.76%  "

# Row 31
$ws.Range("E31").Value = "  +0.00%  "

# Row 32
$ws.Range("D32").Value = "1.47"
$ws.Range("E32").Value = "  +1.38%  "

# Row 33
$ws.Range("D33").Value = "25.53"
$ws.Range("E33").Value = "  +4.79%  "

# Row 34
$ws.Range("D34").Value = "3.557.89"
$ws.Range("E34").Value = "  +4.85%  "

# Row 35
$ws.Range("D35").Value = "0.154"
$ws.Range("E35").Value = "  -3.79%  "

# Row 36
$ws.Range("E36").Value = "  +0.07%  "

# Row 37
$ws.Range("E37").Value = "  +4.40%  "

# Row 38
$ws.Range("D38").Value = "7.91"
$ws.Range("E38").Value = "  +5.80%  "

# Row 39
$ws.Range("D39").Value = "5.61"
$ws.Range("E39").Value = "  +2.02%  "

# Row 40
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  +0.06%  "

# Row 41
$ws.Range("D41").Value = "173.48"
$ws.Range("E41").Value = "  +3.30%  "

# Row 42
$ws.Range("E42").Value = "  +0.77%  "

# Row 43
$ws.Range("D43").Value = "5.21"
$ws.Range("E43").Value = "  +4.01%  "

# Row 44
$ws.Range("D44").Value = "0.895"
$ws.Range("E44").Value = "  +3.01%  "

# Row 45
$ws.Range("E45").Value = "  +0.96%  "

# Row 46
$ws.Range("D46").Value = "46.13"
$ws.Range("E46").Value = "  +1.81%  "

# Row 47
$ws.Range("D47").Value = "1.20"
$ws.Range("E47").Value = "  +1.95%  "

# Row 48
$ws.Range("D48").Value = "25.82"
$ws.Range("E48").Value = "  -1.92%  "

# Row 49
$ws.Range("B49").Value = "Cosmos"
$ws.Range("C49").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D49").Value = "7.14"
$ws.Range("E49").Value = "  +1.81%  "

# Row 50
$ws.Range("B50").Value = "dogwifhat"
$ws.Range("C50").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D50").Value = "2.35"
$ws.Range("E50").Value = "  +4.52%  "

# Row 51
$ws.Range("D51").Value = "23.47"
$ws.Range("E51").Value = "  +16.63%  "
